$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1243.4219
$ws.Range("I15").Value = 1243.4219
$ws.Range("K15").Value = 3730.2657
$ws.Range("M15").Value = -3561.2657

$ws.Range("H28").Value = 998.8461
$ws.Range("I28").Value = 609.3158
$ws.Range("K28").Value = 609.3158
$ws.Range("M28").Value = -124.3158

$ws.Range("H106").Value = 280661.06
$ws.Range("I106").Value = 386915.3
$ws.Range("K106").Value = 386915.3
$ws.Range("M106").Value = -386284.3

$ws.Range("H113").Value = 2599.2727
$ws.Range("I113").Value = 2425.9333
$ws.Range("J113").Value = 2970.7144
$ws.Range("K113").Value = 2425.9333
$ws.Range("L113").Value = 2970.7144
$ws.Range("M113").Value = 828.0666999999999
$ws.Range("N113").Value = -9478.714400000001

$ws.Range("H125").Value = 1949.8182
$ws.Range("I125").Value = 2666.6667
$ws.Range("J125").Value = 1089.6
$ws.Range("K125").Value = 24000.0003
$ws.Range("L125").Value = 9806.4
$ws.Range("M125").Value = -21540.0003
$ws.Range("N125").Value = -14726.4

$ws.Range("H137").Value = 4445278.5
$ws.Range("I137").Value = 793.76666
$ws.Range("J137").Value = 13334248
$ws.Range("K137").Value = 2381.29998
$ws.Range("L137").Value = 40002744
$ws.Range("M137").Value = 168.7000200000002
$ws.Range("N137").Value = -40007844

$ws.Range("H138").Value = 2076.2368
$ws.Range("I138").Value = 1094.25
$ws.Range("J138").Value = 3167.3333
$ws.Range("K138").Value = 3282.75
$ws.Range("L138").Value = 9501.999899999999
$ws.Range("M138").Value = 1857.25
$ws.Range("N138").Value = -19781.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5894.0103
$ws.Range("I32").Value = 5764.2
$ws.Range("K32").Value = 5764.2
$ws.Range("M32").Value = -5477.2

$ws.Range("H74").Value = 13516010
$ws.Range("I74").Value = 19232066
$ws.Range("J74").Value = 5329.8184
$ws.Range("K74").Value = 19232066
$ws.Range("L74").Value = 5329.8184
$ws.Range("M74").Value = -19231192
$ws.Range("N74").Value = -7077.8184

$ws.Range("H77").Value = 13516010
$ws.Range("I77").Value = 19232066
$ws.Range("J77").Value = 5329.8184
$ws.Range("K77").Value = 96160330
$ws.Range("L77").Value = 26649.092
$ws.Range("M77").Value = -96155962
$ws.Range("N77").Value = -35385.092

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 785.2105
$ws.Range("I94").Value = 683
$ws.Range("J94").Value = 1006.6667
$ws.Range("K94").Value = 683
$ws.Range("L94").Value = 1006.6667
$ws.Range("M94").Value = -232
$ws.Range("N94").Value = -1908.6667

$ws.Range("H99").Value = 1218.625
$ws.Range("I99").Value = 1235.5714
$ws.Range("J99").Value = 1100
$ws.Range("K99").Value = 1235.5714
$ws.Range("L99").Value = 1100
$ws.Range("M99").Value = 262.4286
$ws.Range("N99").Value = -4096

$ws.Range("H105").Value = 4636.8184
$ws.Range("J105").Value = 5056.1113
$ws.Range("L105").Value = 5056.1113
$ws.Range("N105").Value = -8550.1113

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5851442
$ws.Range("I31").Value = 3960.0476
$ws.Range("K31").Value = 3960.0476
$ws.Range("M31").Value = -3665.0476

$ws.Range("H34").Value = 5851442
$ws.Range("I34").Value = 3960.0476
$ws.Range("K34").Value = 3960.0476
$ws.Range("M34").Value = -3758.0476

$ws.Range("H107").Value = 842.6
$ws.Range("I107").Value = 674.9286
$ws.Range("J107").Value = 1233.8334
$ws.Range("K107").Value = 674.9286
$ws.Range("L107").Value = 1233.8334
$ws.Range("M107").Value = 1245.0714
$ws.Range("N107").Value = -5073.8334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1107.4783
$ws.Range("I5").Value = 605.1429000000001
$ws.Range("K5").Value = 1815.4287
$ws.Range("M5").Value = -1703.4287

$ws.Range("H105").Value = 1950
$ws.Range("J105").Value = 1950
$ws.Range("L105").Value = 5850
$ws.Range("N105").Value = -11092

$ws.Range("H107").Value = 577.3409
$ws.Range("I107").Value = 151.4
$ws.Range("J107").Value = 631.9487
$ws.Range("K107").Value = 454.2
$ws.Range("L107").Value = 1895.8461
$ws.Range("M107").Value = 1465.8
$ws.Range("N107").Value = -5735.846100000001

$ws.Range("H113").Value = 1884.7858
$ws.Range("I113").Value = 485
$ws.Range("J113").Value = 3751.1667
$ws.Range("K113").Value = 1455
$ws.Range("L113").Value = 11253.5001
$ws.Range("M113").Value = 715
$ws.Range("N113").Value = -15593.5001

$ws.Range("H122").Value = 2029.125
$ws.Range("I122").Value = 2226.8
$ws.Range("K122").Value = 20041.2
$ws.Range("M122").Value = -17591.2

$ws.Range("H129").Value = 2447.4119
$ws.Range("I129").Value = 1739
$ws.Range("J129").Value = 2742.5833
$ws.Range("K129").Value = 5217
$ws.Range("L129").Value = 8227.749899999999
$ws.Range("M129").Value = -217
$ws.Range("N129").Value = -18227.7499

$ws.Range("H131").Value = 852.28
$ws.Range("J131").Value = 859.567
$ws.Range("L131").Value = 2578.701
$ws.Range("N131").Value = -12658.701

$ws.Range("H135").Value = 1107.4783
$ws.Range("I135").Value = 605.1429000000001
$ws.Range("K135").Value = 5446.2861
$ws.Range("M135").Value = -2911.2861

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 59933.332
$ws.Range("J138").Value = 59933.332
$ws.Range("L138").Value = 59933.332
$ws.Range("N138").Value = -70213.33199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1975
$ws.Range("I68").Value = 1750
$ws.Range("J68").Value = 2200
$ws.Range("K68").Value = 1750
$ws.Range("L68").Value = 2200
$ws.Range("M68").Value = -1001
$ws.Range("N68").Value = -3698

$ws.Range("H71").Value = 1975
$ws.Range("I71").Value = 1750
$ws.Range("J71").Value = 2200
$ws.Range("K71").Value = 8750
$ws.Range("L71").Value = 11000
$ws.Range("M71").Value = -5006
$ws.Range("N71").Value = -18488

$ws.Range("H136").Value = 9264699
$ws.Range("J136").Value = 11097.5
$ws.Range("L136").Value = 33292.5
$ws.Range("N136").Value = -38392.5

$ws.Range("H139").Value = 60706.43
$ws.Range("J139").Value = 60706.43
$ws.Range("L139").Value = 60706.43
$ws.Range("N139").Value = -70986.42999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 37424.5
$ws.Range("J42").Value = 37424.5
$ws.Range("L42").Value = 37424.5
$ws.Range("N42").Value = -38180.5

$ws.Range("H136").Value = 2137.0625
$ws.Range("I136").Value = 1668.6923
$ws.Range("J136").Value = 4166.6665
$ws.Range("K136").Value = 5006.0769
$ws.Range("L136").Value = 12499.9995
$ws.Range("M136").Value = -2456.0769
$ws.Range("N136").Value = -17599.9995
